$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newC = 50.788948223972682

$dValues = @{
    2  = 0.99940210580825806
    3  = 0.99946188926696777
    4  = 0.99952167272567749
    5  = 0.99958145618438721
    6  = 0.99964123964309692
    7  = 0.99970108270645142
    8  = 0.99976086616516113
    9  = 0.99982064962387085
    10 = 0.99988043308258057
    11 = 0.99994021654129028
    12 = 0.99994617700576782
    13 = 0.99995219707489014
    14 = 0.99995815753936768
    15 = 0.99996411800384521
    16 = 0.99997007846832275
    17 = 0.99997609853744507
    18 = 0.99998205900192261
    19 = 0.99998801946640015
    20 = 0.99999403953552246
    21 = 0.99999463558197021
    22 = 0.99999523162841797
    23 = 0.99999582767486572
    24 = 0.99999642372131348
    25 = 0.99999701976776123
    26 = 0.99999761581420898
    27 = 0.99999821186065674
    28 = 0.99999880790710449
    29 = 0.99999940395355225
    30 = 0.99999946355819702
    31 = 0.9999995231628418
    32 = 0.99999958276748657
    33 = 0.99999964237213135
    34 = 0.99999970197677612
    35 = 0.9999997615814209
    36 = 0.99999982118606567
    37 = 0.99999988079071045
    38 = 0.99999994039535522
}

for ($row = 2; $row -le 38; $row++) {
    $ws.Cells.Item($row, 3).Value = $newC
    $ws.Cells.Item($row, 4).Value = $dValues[$row]
}
